$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.336.53"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.695.75"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'684.25"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("D6").Value = "'162.65"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("D7").Value = "3.694.54"
$ws.Range("E7").Value = "  -2.64%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "  -6.75%  "
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "'0.450"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  -6.05%  "
$ws.Range("D14").Value = "'33.64"
$ws.Range("E14").Value = "  -6.70%  "
$ws.Range("D15").Value = "4.315.98"
$ws.Range("E15").Value = "  -2.80%  "
$ws.Range("D16").Value = "3.685.47"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "69.437.95"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "'16.36"
$ws.Range("E19").Value = "  -5.78%  "
$ws.Range("E20").Value = "  -7.08%  "
$ws.Range("D21").Value = "'483.77"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "'9.80"
$ws.Range("E22").Value = "  -7.54%  "
$ws.Range("D23").Value = "'0.669"
$ws.Range("E23").Value = "  -8.05%  "
$ws.Range("D24").Value = "'80.08"
$ws.Range("E24").Value = "  -5.73%  "
$ws.Range("D25").Value = "3.836.01"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").Value = "'0.0000129"
$ws.Range("E26").Value = "  -10.70%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'11.52"
$ws.Range("E27").Value = "  -4.65%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -7.94%  "
$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "  -9.68%  "
$ws.Range("E31").Value = "  -10.47%  "
$ws.Range("E32").Value = "  -4.35%  "
$ws.Range("D33").Value = "'6.80"
$ws.Range("E33").Value = "  -6.91%  "
$ws.Range("D34").Value = "'27.11"
$ws.Range("E34").Value = "  -6.59%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'0.166"
$ws.Range("E36").Value = "  -4.53%  "
$ws.Range("D37").Value = "3.654.55"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D38").Value = "'8.59"
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'0.0943"
$ws.Range("E40").Value = "  -6.70%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'2.18"
$ws.Range("E42").Value = "  -5.64%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.963"
$ws.Range("E43").Value = "  -7.09%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'157.95"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("D46").Value = "'48.17"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "'2.84"
$ws.Range("E47").Value = "  -13.24%  "
$ws.Range("D48").Value = "'0.000280"
$ws.Range("E48").Value = "  -12.22%  "
$ws.Range("D49").Value = "'389.79"
$ws.Range("E49").Value = "  -8.16%  "
$ws.Range("D50").Value = "'8.11"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("E51").Value = "  -5.00%  "
